$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 815
$ws1.Range("F5").Value = 973
$ws1.Range("F6").Value = 2326

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 815
$ws4.Range("F7").Value = 973
$ws4.Range("F8").Value = 2326
